$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 239
    $ws.Range("F3").Value = 1062
    $ws.Range("F4").Value = 533
    $ws.Range("G4").Value = 38
    $ws.Range("F5").Value = 13970
    $ws.Range("G7").Value = "已售罄"
    $ws.Range("F8").Value = 226
    $ws.Range("F9").Value = 1799
    $ws.Range("F18").Value = 14044
    $ws.Range("F20").Value = 636
    $ws.Range("F23").Value = 8317
    $ws.Range("F31").Value = 13

    if ($sheetName -eq "展览") {
        $ws.Range("F45").Value = 5118
    } else {
        $ws.Range("F47").Value = 5118
    }
}
